$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 5 (old rows 5-8 shift down to 6-9).
$ws.Rows.Item(5).Insert()

# Fill the newly inserted row 5 with the new job listing.
$ws.Range("B5").Value = "IB報酬を得るための高性能EA開発依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5388547"
$ws.Range("G5").Value = 68
$ws.Range("H5").Value = "◆開発"

# NOTE: the original workbook's hyperlink list is never reconciled against
# row shifts (row 5's hyperlink entry, now stale, is left exactly as-is by
# the real edit being replicated here) -- only a brand-new hyperlink for the
# row that fell off the bottom of the original range is appended
# (old row 8 -> new row 9).
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5388482")
$ws.Range("F9").Style = "Hyperlink"

# Refresh the scrape timestamp on every data row.
$stamp = "2025-09-07 01:21:02"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $stamp
}
